# Update the result of openmpi on the index = 0 dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts E:G -> F:H)
$ws.Columns.Item(5).Insert()

# New header for the inserted column E: "communication time"
$ws.Range("E1").Value = "communication time"

# Fill in the OpenMPI row (row 5) with its results
$ws.Range("B5").Value = 0.150356
$ws.Range("C5").Value = 0.074899
$ws.Range("D5").Value = 0.002739
$ws.Range("E5").Value = 0.000017
$ws.Range("G5").Value = 0.199686
$ws.Range("H5").Value = 0.75407

# Update the selection to G5, matching the author's final cursor position
$ws.Range("G5").Select()

# Set explicit (best-fit style) column width for column F (initCUDA time in train)
$ws.Columns.Item(6).ColumnWidth = 8.43
